# Generate Report for Handoff
# Adds two new handoff records (5479a56b... and b9d5cb6b...) to the
# localization status workbook:
#   - Overview sheet: one row per record (md file entries)
#   - zh-cn / de-de sheets: one row per record (xlf file entries)
# Existing "57f1b5db" rows stay, but shift down to make room for the new
# "5479a56b" record which is inserted right after "9e9f81f3" (row 3),
# and "b9d5cb6b" is appended as the final row.

$wb = $excel.ActiveWorkbook

$repo = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob"

function Set-RowValues($ws, $row, $values) {
    # $values is an ordered list of (col, text)
    foreach ($pair in $values) {
        $col = $pair[0]
        $text = $pair[1]
        $ws.Cells.Item($row, $col).Value = $text
    }
}

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Clear any existing hyperlinks so we can re-add them in final order
$wsOv.Cells.Item(2, 2).Hyperlinks.Delete()

# Row 2 stays the same (9e9f81f3), row 3 becomes the new 5479a56b record,
# row 4 becomes the old 57f1b5db record (shifted down), row 5 is the new
# b9d5cb6b record.
Set-RowValues $wsOv 3 @(
    ,(1, "5479a56b-18b2-4eab-bc44-152efe27dd2a.md")
    ,(2, "e2e\5479a56b-18b2-4eab-bc44-152efe27dd2a.md")
    ,(3, ".md")
    ,(4, "")
    ,(5, "Ready for handoff")
    ,(6, "Ready for handoff")
    ,(7, "2016-08-28 14:40:22")
)

Set-RowValues $wsOv 4 @(
    ,(1, "57f1b5db-41bf-4435-b59b-f28d40c12054.md")
    ,(2, "e2e\57f1b5db-41bf-4435-b59b-f28d40c12054.md")
    ,(3, ".md")
    ,(4, "")
    ,(5, "Ready for handoff")
    ,(6, "Ready for handoff")
    ,(7, "2016-08-28 14:39:08")
)

Set-RowValues $wsOv 5 @(
    ,(1, "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md")
    ,(2, "e2e\b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md")
    ,(3, ".md")
    ,(4, "")
    ,(5, "Ready for handoff")
    ,(6, "Ready for handoff")
    ,(7, "2016-08-28 14:40:22")
)

# Re-create the hyperlinks for column B in row order (2,3,4,5)
$wsOv.Hyperlinks.Add($wsOv.Cells.Item(2, 2), "$repo/0208d729dfff7e4089f6cee3f44a9d8f40c64d54/e2e/9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md", "", "", "e2e\9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Cells.Item(3, 2), "$repo/29ab18afe0c0476a09d22c4df95a3e3b78844a77/e2e/5479a56b-18b2-4eab-bc44-152efe27dd2a.md", "", "", "e2e\5479a56b-18b2-4eab-bc44-152efe27dd2a.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Cells.Item(4, 2), "$repo/34fb509076e49b457389259e65d9400024201c2a/e2e/57f1b5db-41bf-4435-b59b-f28d40c12054.md", "", "", "e2e\57f1b5db-41bf-4435-b59b-f28d40c12054.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Cells.Item(5, 2), "$repo/9565b506f503a9138b4ccb2f240843822587010a/e2e/b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md", "", "", "e2e\b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md") | Out-Null

# Grow the Overview table (table3) to the new extents
$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G5"))

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(3, 1).Hyperlinks.Delete()

Set-RowValues $wsZh 3 @(
    ,(1, "5479a56b-18b2-4eab-bc44-152efe27dd2a.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "5479a56b-18b2-4eab-bc44-152efe27dd2a.29ab18afe0c0476a09d22c4df95a3e3b78844a77.zh-cn.xlf")
    ,(8, "2016-08-28 14:40:18")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

Set-RowValues $wsZh 4 @(
    ,(1, "57f1b5db-41bf-4435-b59b-f28d40c12054.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "57f1b5db-41bf-4435-b59b-f28d40c12054.4d872bec399f7bf0603d4e693415e62ade1d31bb.zh-cn.xlf")
    ,(8, "2016-08-28 14:39:00")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

Set-RowValues $wsZh 5 @(
    ,(1, "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.9565b506f503a9138b4ccb2f240843822587010a.zh-cn.xlf")
    ,(8, "2016-08-28 14:40:18")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 1), "$repo/0208d729dfff7e4089f6cee3f44a9d8f40c64d54/e2e/9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md", "", "", "9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e27563d0d3eb4c4811d3ee1be5d6a9197d91a57e/e2e/9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md", "", "", "9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 1), "$repo/29ab18afe0c0476a09d22c4df95a3e3b78844a77/e2e/5479a56b-18b2-4eab-bc44-152efe27dd2a.md", "", "", "5479a56b-18b2-4eab-bc44-152efe27dd2a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 1), "$repo/34fb509076e49b457389259e65d9400024201c2a/e2e/57f1b5db-41bf-4435-b59b-f28d40c12054.md", "", "", "57f1b5db-41bf-4435-b59b-f28d40c12054.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 1), "$repo/9565b506f503a9138b4ccb2f240843822587010a/e2e/b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md", "", "", "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md") | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(3, 1).Hyperlinks.Delete()

Set-RowValues $wsDe 3 @(
    ,(1, "5479a56b-18b2-4eab-bc44-152efe27dd2a.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "5479a56b-18b2-4eab-bc44-152efe27dd2a.29ab18afe0c0476a09d22c4df95a3e3b78844a77.de-de.xlf")
    ,(8, "2016-08-28 14:40:22")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

Set-RowValues $wsDe 4 @(
    ,(1, "57f1b5db-41bf-4435-b59b-f28d40c12054.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "57f1b5db-41bf-4435-b59b-f28d40c12054.4d872bec399f7bf0603d4e693415e62ade1d31bb.de-de.xlf")
    ,(8, "2016-08-28 14:39:08")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

Set-RowValues $wsDe 5 @(
    ,(1, "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md")
    ,(2, ".md")
    ,(3, "Ready for handoff")
    ,(4, "e2e")
    ,(5, "ht")
    ,(6, "False")
    ,(7, "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.9565b506f503a9138b4ccb2f240843822587010a.de-de.xlf")
    ,(8, "2016-08-28 14:40:22")
    ,(9, "")
    ,(10, "")
    ,(11, "0001-01-01 00:00:00")
    ,(12, "")
    ,(13, "True")
    ,(14, "")
    ,(15, "False")
    ,(16, "")
)

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 1), "$repo/0208d729dfff7e4089f6cee3f44a9d8f40c64d54/e2e/9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md", "", "", "9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2bafd66c7f0a5389a77a3c13af0ec334b32ccb12/e2e/9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md", "", "", "9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 1), "$repo/29ab18afe0c0476a09d22c4df95a3e3b78844a77/e2e/5479a56b-18b2-4eab-bc44-152efe27dd2a.md", "", "", "5479a56b-18b2-4eab-bc44-152efe27dd2a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 1), "$repo/34fb509076e49b457389259e65d9400024201c2a/e2e/57f1b5db-41bf-4435-b59b-f28d40c12054.md", "", "", "57f1b5db-41bf-4435-b59b-f28d40c12054.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 1), "$repo/9565b506f503a9138b4ccb2f240843822587010a/e2e/b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md", "", "", "b9d5cb6b-ddc4-4e75-96f4-3cb3b24c5a11.md") | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

Write-Host "Done"
